$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text temporarily so numeric-looking price strings
# (e.g. "1.00", "230.12") are stored as text, matching the source data,
# instead of being auto-converted to numbers.
$priceRange = $ws.Range('D2:D51')
$priceRange.NumberFormat = '@'

$ws.Range('D2').Value = '43.182.70'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.237.09'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '230.12'
$ws.Range('E5').Value = '  -1.85%  '
$ws.Range('D6').Value = '0.639'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').Value = '62.60'
$ws.Range('E7').Value = '  -4.96%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.439'
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('D10').Value = '0.0955'
$ws.Range('E10').Value = '  -8.03%  '
$ws.Range('D11').Value = '56.58'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('D12').Value = '26.80'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').Value = '2.565.23'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').Value = '15.33'
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').Value = '5.99'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = '0.821'
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = '2.237.14'
$ws.Range('E18').Value = '  -2.57%  '
$ws.Range('D19').Value = '43.063.42'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('D20').Value = '0.0₃0961'
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('D21').Value = '72.72'
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').Value = '6.05'
$ws.Range('E22').Value = '  -2.04%  '
$ws.Range('D23').Value = '244.62'
$ws.Range('E23').Value = '  -6.58%  '
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  +30.76%  '
$ws.Range('E26').Value = '  -5.02%  '
$ws.Range('D27').Value = '2.20'
$ws.Range('E27').Value = '  -5.44%  '
$ws.Range('D28').Value = '9.70'
$ws.Range('E28').Value = '  -4.88%  '
$ws.Range('D29').Value = '171.33'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('D30').Value = '21.48'
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('D31').Value = '0.129'
$ws.Range('E31').Value = '  -5.85%  '
$ws.Range('D32').Value = '1.40'
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('D34').Value = '4.83'
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('D35').Value = '0.0668'
$ws.Range('E35').Value = '  -3.04%  '
$ws.Range('D36').Value = '4.83'
$ws.Range('E36').Value = '  -5.49%  '
$ws.Range('D37').Value = '3.59'
$ws.Range('E37').Value = '  -6.00%  '
$ws.Range('D38').Value = '6.28'
$ws.Range('E38').Value = '  -7.89%  '
$ws.Range('D39').Value = '2.26'
$ws.Range('E39').Value = '  -4.73%  '
$ws.Range('D40').Value = '0.0249'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '8.50'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('D43').Value = '4.47'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').Value = '16.86'
$ws.Range('E44').Value = '  -3.83%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '0.0934'
$ws.Range('E45').Value = '  -4.78%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '95.53'
$ws.Range('E46').Value = '  -3.40%  '
$ws.Range('B47').Value = 'TerraClassic'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D47').Value = '0.000208'
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = '1.17'
$ws.Range('E48').Value = '  -3.43%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.442.36'
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '2.26'
$ws.Range('E50').Value = '  -4.91%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = '2.74'
$ws.Range('E51').Value = '  +1.08%  '

# Restore default (General) style so no stray formatting is left behind
# on the price column (values remain text since they were already entered).
$priceRange.Style = 'Normal'
